# "added role to users list"
#
# Populate the 従業員 (Employee) sheet with four employee rows (2-5) and
# seed the default values shown on the 部署 (Department) / 役職
# (Position) lookup sheets.
#
# Columns on 従業員: A=社員番号 B=姓 C=名 D=表示名 E=権限 F=パスワード
#                    G=部署 H=役職 I=電話番号 J=メールアドレス

$wb = $excel.ActiveWorkbook

$wsEmployee   = $wb.Worksheets.Item("従業員")
$wsDepartment = $wb.Worksheets.Item("部署")
$wsPosition   = $wb.Worksheets.Item("役職")

function Set-TextValue($range, $value) {
    # Force a literal value to be stored as text even when it looks like a
    # number (e.g. phone numbers), without disturbing the cell's existing
    # number format / style.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
}

# --- Row 2 -------------------------------------------------------------
$wsEmployee.Range("A2").Value = "DK0002"
$wsEmployee.Range("B2").Value = "dsadasdsad"
$wsEmployee.Range("C2").Value = "sasadsa"
$wsEmployee.Range("D2").Value = "dsas"
$wsEmployee.Range("E2").Value = "EMPLOYEE"
$wsEmployee.Range("F2").Value = ""
$wsEmployee.Range("G2").Value = "営業部"
$wsEmployee.Range("H2").Value = "開発部"
$wsEmployee.Range("I2").Value = ""
$wsEmployee.Range("J2").Value = "justcheckout08@gmail.com"

# --- Row 3 -------------------------------------------------------------
$wsEmployee.Range("A3").Value = "DK0001"
$wsEmployee.Range("B3").Value = "MAMBETLEPESOV"
$wsEmployee.Range("C3").Value = "RINAT"
$wsEmployee.Range("D3").Value = "Alik"
$wsEmployee.Range("E3").Value = "EMPLOYEE"
$wsEmployee.Range("F3").Value = ""
$wsEmployee.Range("G3").Value = "営業部"
$wsEmployee.Range("H3").Value = "開発部"
$wsEmployee.Range("I3").Value = ""
$wsEmployee.Range("J3").Value = "rinatmambetlepesov@gmail.com"

# --- Row 4 -------------------------------------------------------------
$wsEmployee.Range("A4").Value = "Cat"
$wsEmployee.Range("B4").Value = "MAMBETLEPESOV"
$wsEmployee.Range("C4").Value = "RINAT"
$wsEmployee.Range("D4").Value = "Alikaa"
$wsEmployee.Range("E4").Value = "EMPLOYEE"
$wsEmployee.Range("F4").Value = ""
$wsEmployee.Range("G4").Value = "営業部"
$wsEmployee.Range("H4").Value = "開発部"
Set-TextValue $wsEmployee.Range("I4") "321321"
$wsEmployee.Range("J4").Value = "rinatmambetlepeso@gmail.com"

# --- Row 5 -------------------------------------------------------------
$wsEmployee.Range("A5").Value = "Admin02"
$wsEmployee.Range("B5").Value = "MAMBETLEPESOV"
$wsEmployee.Range("C5").Value = "RINAT"
$wsEmployee.Range("D5").Value = ""
$wsEmployee.Range("E5").Value = "EMPLOYEE"
$wsEmployee.Range("F5").Value = ""
$wsEmployee.Range("G5").Value = "営業部"
$wsEmployee.Range("H5").Value = "開発部"
Set-TextValue $wsEmployee.Range("I5") "939207055"
$wsEmployee.Range("J5").Value = "sarvarbekfozilov59@gmail.com"

# --- 部署 (Department) lookup sheet ------------------------------------
$wsDepartment.Range("A2").Value = "営業部"

# --- 役職 (Position) lookup sheet ---------------------------------------
$wsPosition.Range("A2").Value = "開発部"
